$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.062.62'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.842.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7009'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.33%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3033'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07438'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.02%  '
$ws.Range("E10").Value = '  -6.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08106'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7223'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.241'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.814.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '88.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.875.91'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.780'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '239.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007650'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9993'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.045.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.528'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.63%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1459'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.77%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.928'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.94'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.14%  '
$ws.Range("E29").Value = '  -5.30%  '
$ws.Range("E30").Value = '  -7.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.471'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.11%  '
$ws.Range("E32").Value = '  -3.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.010'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05156'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.180'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7057'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.017'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.37%  '
$ws.Range("E38").Value = '  -2.18%  '
$ws.Range("E39").Value = '  -5.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.668'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8938'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.917'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4274'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.058.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.748'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.29%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.197'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.65%  '
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.038'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.58%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.966.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.96%  '
